# Update Name of Algo
# Apply updated imputed values to the RandomForest result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4"  = -12.45329999999999
    "B9"  = 6.258399999999994
    "C9"  = -10.9742
    "C11" = -12.9673
    "B18" = 6.216899999999998
    "B20" = 9.059299999999995
    "C23" = -12.1685
    "C24" = -13.17839999999999
    "C26" = -12.8794
    "B27" = 6.434100000000004
    "C34" = -12.0642
    "B35" = 8.832699999999997
    "C35" = -13.20190000000001
    "C48" = -11.38559999999999
    "C49" = -13.7757
    "C52" = -10.93669999999999
    "C66" = -11.12540000000001
    "C67" = -11.2902
    "B69" = 5.488999999999993
    "B76" = 5.5394
    "B78" = 10.356
    "C78" = -13.40150000000001
    "C80" = -13.11880000000001
    "B82" = 5.391099999999998
    "B83" = 5.235799999999998
    "B93" = 5.4265
    "C99" = -12.9414
    "C104" = -12.80200000000001
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
